$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "spawn_npc_inner" row (row 7) entirely, shifting the rows below it up.
$ws.Rows("7").Delete()

# Move selection to match the saved selection state of the edited file.
$ws.Range("F20").Select()
